$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 141.58241910315
$ws.Range("C3").Value = 10.60823564375479
$ws.Range("C4").Value = 6.689209915980248
$ws.Range("C5").Value = 15.55895387579201
$ws.Range("C6").Value = 38.70213629730312
$ws.Range("C7").Value = 11.97263140828689
$ws.Range("C8").Value = 9.011884789316531
$ws.Range("C9").Value = 27.88693735784083
$ws.Range("C10").Value = 43.67081738142062
$ws.Range("C11").Value = 9.648394547813192
$ws.Range("C12").Value = 3.919806721343242
$ws.Range("C13").Value = 7.210132626308191
$ws.Range("C14").Value = 1.703346973351188
$ws.Range("C15").Value = 1.654144378522611
$ws.Range("C16").Value = 20.94702850630546
$ws.Range("C17").Value = 19.33583877406176
$ws.Range("C18").Value = 9.960791975296218
$ws.Range("C19").Value = 1.091829009053168
$ws.Range("C20").Value = 29.51687093573351
$ws.Range("C21").Value = 76.40147685238573
$ws.Range("C22").Value = 14.35934775425719
$ws.Range("C23").Value = 0.2319550899061452
$ws.Range("C24").Value = 2.562439898929504
$ws.Range("C25").Value = 26.97317488245299
$ws.Range("C26").Value = 6.894611224550338
$ws.Range("C27").Value = 0.5146747617782819
$ws.Range("C28").Value = 10.74490951827861
$ws.Range("C29").Value = 24.58333456220786
$ws.Range("C30").Value = 9.985783769494859
$ws.Range("C31").Value = 12.93169151065977
$ws.Range("C32").Value = 3.576950544680624
$ws.Range("C33").Value = 1.83377289932535
$ws.Range("C34").Value = 5.022569640358316
$ws.Range("C35").Value = 2.409365159462822
$ws.Range("C36").Value = 85.98582992756491
$ws.Range("C37").Value = 8.138733979501479
$ws.Range("C38").Value = 25.37448104730862
$ws.Range("C39").Value = 5.09520204224812
$ws.Range("C40").Value = 3.530871924126878
$ws.Range("C41").Value = 12.67162065228016
$ws.Range("C42").Value = 0.8262911956925982
$ws.Range("C43").Value = 5.799658241222338
$ws.Range("C44").Value = 242.879842700693
